$d = $word.ActiveDocument

# ----------------------------------------------------------------------------
# Helper: replace the whole text of a paragraph while preserving any leading
# empty run (e.g. <w:r/>) that precedes the text run. A plain Find/Replace can
# merge an empty run into the text run when both runs carry no character
# formatting, so for those paragraphs we rebuild just the text run via
# Range.InsertXML (which only touches the exact range supplied).
# ----------------------------------------------------------------------------
function Set-ParagraphText($paragraphIndex, $newText) {
    $paragraph = $d.Paragraphs($paragraphIndex)
    $start = $paragraph.Range.Start
    $end = $paragraph.Range.End
    $r = $d.Range($start, $end)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>' + $newText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

function Find-ParagraphIndex($text) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        if ($d.Paragraphs($i).Range.Text -eq $text) {
            return $i
        }
    }
    return -1
}

# ----------------------------------------------------------------------------
# Title (Heading1) and its bold repeat near the end of the document share the
# exact same original text, so a single document-wide Find/Replace-all fixes
# both occurrences at once. The bold run keeps its own formatting (and the
# empty run ahead of it survives because that empty run's formatting differs
# from the bold run's formatting).
# ----------------------------------------------------------------------------
$d.Content.Find.Execute("Play Free Miner Donkey Trouble Slot - Review and Gameplay", $true, $true, $false, $false, $false, $true, 1, $false, "Play Miner Donkey Trouble for Free", 2)

# ----------------------------------------------------------------------------
# "What we like" bullet list (each paragraph begins with an empty <w:r/>)
# ----------------------------------------------------------------------------
$i = Find-ParagraphIndex "Attractive cartoon-style graphics and animations`r"
Set-ParagraphText $i "Exciting gameplay with unique theme"

$i = Find-ParagraphIndex "Exciting features like exploding symbols and free spins`r"
Set-ParagraphText $i "Attractive graphics and animations"

$i = Find-ParagraphIndex "Selectable betting amounts for players of all levels`r"
Set-ParagraphText $i "Wide range of betting options"

$i = Find-ParagraphIndex "High-stakes gambling with high volatility`r"
Set-ParagraphText $i "Generous payouts with exciting features"

# ----------------------------------------------------------------------------
# "What we don't like" bullet list
# ----------------------------------------------------------------------------
$i = Find-ParagraphIndex "Slightly lower RTP value compared to other slot games`r"
Set-ParagraphText $i "Lower RTP value compared to other slots"

$i = Find-ParagraphIndex "Winning combinations start from 5 identical symbols, which may be difficult to achieve`r"
Set-ParagraphText $i "High volatility may not be suitable for all players"

# ----------------------------------------------------------------------------
# Italic meta description near the very end of the document.
# ----------------------------------------------------------------------------
$d.Content.Find.Execute("Discover the world of mining with Miner Donkey Trouble slot game. Exciting features, selectable betting, and high volatility. Play for free and enjoy!", $true, $true, $false, $false, $false, $true, 1, $false, "Read our review of Miner Donkey Trouble and play this exciting slot game for free.", 2)
